$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): E1/F1/G1 ---
$ws.Range("E1").Value = "Đơn vị mặc định"
$ws.Range("F1").Value = "Đơn vị mua"
$ws.Range("G1").Value = "Mức tồn tối thiểu"

# Match header styling of the existing header row (bold, like A1:C1/D1)
$ws.Range("E1:G1").Font.Bold = $true

# --- New data cells (row 2): E2/F2/G2 ---
$ws.Range("E2").Value = "Cái"
$ws.Range("F2").Value = "Cái"
$ws.Range("G2").Value = 100

# --- Column widths for the newly added columns (best-fit sized) ---
$ws.Columns.Item(5).ColumnWidth = 14.76
$ws.Columns.Item(6).ColumnWidth = 10.42
$ws.Columns.Item(7).ColumnWidth = 15.75

# --- Selection moves to H2 after the edit ---
$ws.Range("H2").Select()
